# Natmi following Dr Hou advice
# Recompute the Fzd8 -> Ckap4 LR-pair table for all Sending/Target cluster
# combinations (ECs, FAPs, sCs) instead of only the previous subset, using
# the updated (Dr Hou advised) statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-D: Sending cluster, Ligand symbol, Receptor symbol, Target cluster
$labels = @(
  @("ECs",  "Fzd8", "Ckap4", "ECs"),
  @("ECs",  "Fzd8", "Ckap4", "FAPs"),
  @("ECs",  "Fzd8", "Ckap4", "sCs"),
  @("FAPs", "Fzd8", "Ckap4", "ECs"),
  @("FAPs", "Fzd8", "Ckap4", "FAPs"),
  @("FAPs", "Fzd8", "Ckap4", "sCs"),
  @("sCs",  "Fzd8", "Ckap4", "ECs"),
  @("sCs",  "Fzd8", "Ckap4", "FAPs"),
  @("sCs",  "Fzd8", "Ckap4", "sCs")
)

# Columns E-T (16 numeric metrics per row)
$values = @(
  @(2, 0.6666666666666666, 2.558821666666667, 7.676465, 0.2156728774407755, 0.2156728774407755, 3, 1, 10.36378, 31.09134, 0.2440702244474787, 0.2440702244474788, 26.51906481256667, 238.6715833131, 0.05263932760420364, 0.05263932760420365),
  @(2, 0.6666666666666666, 2.558821666666667, 7.676465, 0.2156728774407755, 0.2156728774407755, 3, 1, 19.886525, 59.659575, 0.4683338145184863, 0.4683338145184862, 50.88607104470834, 457.974639402375, 0.1010069013800163, 0.1010069013800163),
  @(2, 0.6666666666666666, 2.558821666666667, 7.676465, 0.2156728774407755, 0.2156728774407755, 3, 1, 12.21198233333333, 36.635947, 0.287595961034035, 0.2875959610340351, 31.24828498748389, 281.234564887355, 0.06202664845655546, 0.06202664845655548),
  @(3, 1, 6.453984666666667, 19.361954, 0.5439806384912759, 0.5439806384912759, 3, 1, 10.36378, 31.09134, 0.2440702244474787, 0.2440702244474788, 66.88767720870668, 601.9890948783601, 0.1327694765316485, 0.1327694765316485),
  @(3, 1, 6.453984666666667, 19.361954, 0.5439806384912759, 0.5439806384912759, 3, 1, 19.886525, 59.659575, 0.4683338145184863, 0.4683338145184862, 128.3473274232834, 1155.12594680955, 0.2547645274488209, 0.2547645274488209),
  @(3, 1, 6.453984666666667, 19.361954, 0.5439806384912759, 0.5439806384912759, 3, 1, 12.21198233333333, 36.635947, 0.287595961034035, 0.2875959610340351, 78.81594672893756, 709.343520560438, 0.1564466345108065, 0.1564466345108065),
  @(3, 1, 2.851558333333334, 8.554675000000001, 0.2403464840679487, 0.2403464840679487, 3, 1, 10.36378, 31.09134, 0.2440702244474787, 0.2440702244474788, 29.55292322383334, 265.9763090145001, 0.05866142031162661, 0.05866142031162662),
  @(3, 1, 2.851558333333334, 8.554675000000001, 0.2403464840679487, 0.2403464840679487, 3, 1, 19.886525, 59.659575, 0.4683338145184863, 0.4683338145184862, 56.70758608479168, 510.3682747631251, 0.112562385689649, 0.112562385689649),
  @(3, 1, 2.851558333333334, 8.554675000000001, 0.2403464840679487, 0.2403464840679487, 3, 1, 12.21198233333333, 36.635947, 0.287595961034035, 0.2875959610340351, 34.82317998913611, 313.408619902225, 0.06912267806667309, 0.06912267806667312)
)

$startRow = 2
for ($i = 0; $i -lt $labels.Length; $i++) {
  $r = $startRow + $i

  $lbl = $labels[$i]
  for ($c = 0; $c -lt $lbl.Length; $c++) {
    $ws.Cells.Item($r, $c + 1).Value2 = $lbl[$c]
  }

  $vals = $values[$i]
  for ($c = 0; $c -lt $vals.Length; $c++) {
    $ws.Cells.Item($r, $c + 5).Value2 = $vals[$c]
  }
}
